# Update scripts with new TPM values.
# - Adds a new "Resolving-Mac" sending-cluster cohort (rows 14-16).
# - Refreshes the TPM-derived NATMI metrics (columns E:T) for every existing
#   sending/target cluster combination now that the "Gdf9" / "Bmpr1b" shared
#   strings are recomputed against the new MuSCs-adjacent "Resolving-Mac" group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gdf9"
$ws.Range("C2").Value = "Bmpr1b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.612610333333333
$ws.Range("H2").Value = 4.837831
$ws.Range("I2").Value = 0.2199512280440611
$ws.Range("J2").Value = 0.2199512280440611
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.07259900000000001
$ws.Range("N2").Value = 0.217797
$ws.Range("O2").Value = 0.0162094769588191
$ws.Range("P2").Value = 0.0162094769588191
$ws.Range("Q2").Value = 0.1170738975896667
$ws.Range("R2").Value = 1.053665078307
$ws.Range("S2").Value = 0.003565294363044173
$ws.Range("T2").Value = 0.003565294363044174

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gdf9"
$ws.Range("C3").Value = "Bmpr1b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.612610333333333
$ws.Range("H3").Value = 4.837831
$ws.Range("I3").Value = 0.2199512280440611
$ws.Range("J3").Value = 0.2199512280440611
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.187332
$ws.Range("N3").Value = 12.561996
$ws.Range("O3").Value = 0.9349228167457665
$ws.Range("P3").Value = 0.9349228167457664
$ws.Range("Q3").Value = 6.752534852297334
$ws.Range("R3").Value = 60.77281367067599
$ws.Range("S3").Value = 0.205637421669644
$ws.Range("T3").Value = 0.205637421669644

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gdf9"
$ws.Range("C4").Value = "Bmpr1b"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.612610333333333
$ws.Range("H4").Value = 4.837831
$ws.Range("I4").Value = 0.2199512280440611
$ws.Range("J4").Value = 0.2199512280440611
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2188686666666667
$ws.Range("N4").Value = 0.656606
$ws.Range("O4").Value = 0.04886770629541442
$ws.Range("P4").Value = 0.04886770629541441
$ws.Range("Q4").Value = 0.3529498735095556
$ws.Range("R4").Value = 3.176548861586
$ws.Range("S4").Value = 0.0107485120113729
$ws.Range("T4").Value = 0.0107485120113729

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gdf9"
$ws.Range("C5").Value = "Bmpr1b"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 4.376645666666667
$ws.Range("H5").Value = 13.129937
$ws.Range("I5").Value = 0.5969505274762915
$ws.Range("J5").Value = 0.5969505274762916
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.07259900000000001
$ws.Range("N5").Value = 0.217797
$ws.Range("O5").Value = 0.0162094769588191
$ws.Range("P5").Value = 0.0162094769588191
$ws.Range("Q5").Value = 0.3177400987543334
$ws.Range("R5").Value = 2.859660888789
$ws.Range("S5").Value = 0.009676255820681856
$ws.Range("T5").Value = 0.009676255820681857

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gdf9"
$ws.Range("C6").Value = "Bmpr1b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 4.376645666666667
$ws.Range("H6").Value = 13.129937
$ws.Range("I6").Value = 0.5969505274762915
$ws.Range("J6").Value = 0.5969505274762916
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.187332
$ws.Range("N6").Value = 12.561996
$ws.Range("O6").Value = 0.9349228167457665
$ws.Range("P6").Value = 0.9349228167457664
$ws.Range("Q6").Value = 18.32646845269467
$ws.Range("R6").Value = 164.938216074252
$ws.Range("S6").Value = 0.5581026686060055
$ws.Range("T6").Value = 0.5581026686060055

# Row 7: FAPs -> MuSCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Gdf9"
$ws.Range("C7").Value = "Bmpr1b"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 4.376645666666667
$ws.Range("H7").Value = 13.129937
$ws.Range("I7").Value = 0.5969505274762915
$ws.Range("J7").Value = 0.5969505274762916
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2188686666666667
$ws.Range("N7").Value = 0.656606
$ws.Range("O7").Value = 0.04886770629541442
$ws.Range("P7").Value = 0.04886770629541441
$ws.Range("Q7").Value = 0.9579106015357779
$ws.Range("R7").Value = 8.621195413822
$ws.Range("S7").Value = 0.02917160304960413
$ws.Range("T7").Value = 0.02917160304960413

# Row 8: Inflammatory-Mac -> ECs
$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("B8").Value = "Gdf9"
$ws.Range("C8").Value = "Bmpr1b"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.1451676666666667
$ws.Range("H8").Value = 0.435503
$ws.Range("I8").Value = 0.01980007562622025
$ws.Range("J8").Value = 0.01980007562622025
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.07259900000000001
$ws.Range("N8").Value = 0.217797
$ws.Range("O8").Value = 0.0162094769588191
$ws.Range("P8").Value = 0.0162094769588191
$ws.Range("Q8").Value = 0.01053902743233334
$ws.Range("R8").Value = 0.094851246891
$ws.Range("S8").Value = 0.0003209488696460928
$ws.Range("T8").Value = 0.0003209488696460928

# Row 9: Inflammatory-Mac -> FAPs
$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("B9").Value = "Gdf9"
$ws.Range("C9").Value = "Bmpr1b"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.1451676666666667
$ws.Range("H9").Value = 0.435503
$ws.Range("I9").Value = 0.01980007562622025
$ws.Range("J9").Value = 0.01980007562622025
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.187332
$ws.Range("N9").Value = 12.561996
$ws.Range("O9").Value = 0.9349228167457665
$ws.Range("P9").Value = 0.9349228167457664
$ws.Range("Q9").Value = 0.6078652159986667
$ws.Range("R9").Value = 5.470786943988
$ws.Range("S9").Value = 0.01851154247624503
$ws.Range("T9").Value = 0.01851154247624503

# Row 10: Inflammatory-Mac -> MuSCs
$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("B10").Value = "Gdf9"
$ws.Range("C10").Value = "Bmpr1b"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.1451676666666667
$ws.Range("H10").Value = 0.435503
$ws.Range("I10").Value = 0.01980007562622025
$ws.Range("J10").Value = 0.01980007562622025
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2188686666666667
$ws.Range("N10").Value = 0.656606
$ws.Range("O10").Value = 0.04886770629541442
$ws.Range("P10").Value = 0.04886770629541441
$ws.Range("Q10").Value = 0.03177265364644444
$ws.Range("R10").Value = 0.285953882818
$ws.Range("S10").Value = 0.0009675842803291247
$ws.Range("T10").Value = 0.0009675842803291248

# Row 11: MuSCs -> ECs
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Gdf9"
$ws.Range("C11").Value = "Bmpr1b"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.6315603333333334
$ws.Range("H11").Value = 1.894681
$ws.Range("I11").Value = 0.08614137465772362
$ws.Range("J11").Value = 0.08614137465772362
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.07259900000000001
$ws.Range("N11").Value = 0.217797
$ws.Range("O11").Value = 0.0162094769588191
$ws.Range("P11").Value = 0.0162094769588191
$ws.Range("Q11").Value = 0.04585064863966667
$ws.Range("R11").Value = 0.412655837757
$ws.Range("S11").Value = 0.001396306627715375
$ws.Range("T11").Value = 0.001396306627715375

# Row 12: MuSCs -> FAPs
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Gdf9"
$ws.Range("C12").Value = "Bmpr1b"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.6315603333333334
$ws.Range("H12").Value = 1.894681
$ws.Range("I12").Value = 0.08614137465772362
$ws.Range("J12").Value = 0.08614137465772362
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 4.187332
$ws.Range("N12").Value = 12.561996
$ws.Range("O12").Value = 0.9349228167457665
$ws.Range("P12").Value = 0.9349228167457664
$ws.Range("Q12").Value = 2.644552793697334
$ws.Range("R12").Value = 23.800975143276
$ws.Range("S12").Value = 0.08053553663335135
$ws.Range("T12").Value = 0.08053553663335135

# Row 13: MuSCs -> MuSCs
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Gdf9"
$ws.Range("C13").Value = "Bmpr1b"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.6315603333333334
$ws.Range("H13").Value = 1.894681
$ws.Range("I13").Value = 0.08614137465772362
$ws.Range("J13").Value = 0.08614137465772362
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.2188686666666667
$ws.Range("N13").Value = 0.656606
$ws.Range("O13").Value = 0.04886770629541442
$ws.Range("P13").Value = 0.04886770629541441
$ws.Range("Q13").Value = 0.1382287680762223
$ws.Range("R13").Value = 1.244058912686
$ws.Range("S13").Value = 0.004209531396656893
$ws.Range("T13").Value = 0.004209531396656892

# Row 14: Resolving-Mac -> ECs
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Gdf9"
$ws.Range("C14").Value = "Bmpr1b"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.5656883333333332
$ws.Range("H14").Value = 1.697065
$ws.Range("I14").Value = 0.0771567941957035
$ws.Range("J14").Value = 0.07715679419570351
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.07259900000000001
$ws.Range("N14").Value = 0.217797
$ws.Range("O14").Value = 0.0162094769588191
$ws.Range("P14").Value = 0.0162094769588191
$ws.Range("Q14").Value = 0.04106840731166667
$ws.Range("R14").Value = 0.369615665805
$ws.Range("S14").Value = 0.001250671277731603
$ws.Range("T14").Value = 0.001250671277731603

# Row 15: Resolving-Mac -> FAPs
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Gdf9"
$ws.Range("C15").Value = "Bmpr1b"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.5656883333333332
$ws.Range("H15").Value = 1.697065
$ws.Range("I15").Value = 0.0771567941957035
$ws.Range("J15").Value = 0.07715679419570351
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 4.187332
$ws.Range("N15").Value = 12.561996
$ws.Range("O15").Value = 0.9349228167457665
$ws.Range("P15").Value = 0.9349228167457664
$ws.Range("Q15").Value = 2.368724860193333
$ws.Range("R15").Value = 21.31852374174
$ws.Range("S15").Value = 0.07213564736052051
$ws.Range("T15").Value = 0.07213564736052053

# Row 16: Resolving-Mac -> MuSCs
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Gdf9"
$ws.Range("C16").Value = "Bmpr1b"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.5656883333333332
$ws.Range("H16").Value = 1.697065
$ws.Range("I16").Value = 0.0771567941957035
$ws.Range("J16").Value = 0.07715679419570351
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2188686666666667
$ws.Range("N16").Value = 0.656606
$ws.Range("O16").Value = 0.04886770629541442
$ws.Range("P16").Value = 0.04886770629541441
$ws.Range("Q16").Value = 0.1238114512655555
$ws.Range("R16").Value = 1.11430306139
$ws.Range("S16").Value = 0.003770475557451374
$ws.Range("T16").Value = 0.003770475557451375

